$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for the numeric-looking Price/Volume columns so COM
# auto-detection does not coerce values like "0.9984" or "1.000" into numbers.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2: Bitcoin
$ws.Range("D2").Value = "29.044.47"
$ws.Range("E2").Value = "  -0.51%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "1.829.25"
$ws.Range("E3").Value = "  -0.28%  "

# Row 4: TetherUSD
$ws.Range("D4").Value = "0.9984"

# Row 5: BNB
$ws.Range("D5").Value = "240.61"
$ws.Range("E5").Value = "  -0.50%  "

# Row 6: XRP
$ws.Range("D6").Value = "0.6207"
$ws.Range("E6").Value = "  -6.73%  "

# Row 7: USDC
$ws.Range("E7").Value = "  +0.02%  "

# Row 8: Dogecoin -> OKB
$ws.Range("B8").Value = "OKB"
$ws.Range("C8").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D8").Value = "44.50"
$ws.Range("E8").Value = "  +6.18%  "

# Row 9: OKB -> Dogecoin
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "0.07489"
$ws.Range("E9").Value = "  +0.97%  "

# Row 10: Cardano
$ws.Range("E10").Value = "  -0.85%  "

# Row 11: Solana
$ws.Range("D11").Value = "22.70"
$ws.Range("E11").Value = "  -1.20%  "

# Row 12: TRON
$ws.Range("D12").Value = "0.07624"
$ws.Range("E12").Value = "  -1.71%  "

# Row 13: WrappedEther
$ws.Range("D13").Value = "1.826.94"
$ws.Range("E13").Value = "  -0.95%  "

# Row 14: Polkadot
$ws.Range("D14").Value = "4.950"
$ws.Range("E14").Value = "  -0.86%  "

# Row 15: Polygon
$ws.Range("D15").Value = "0.6627"
$ws.Range("E15").Value = "  -0.87%  "

# Row 16: Litecoin
$ws.Range("D16").Value = "82.07"
$ws.Range("E16").Value = "  -1.08%  "

# Row 17: ShibaInu
$ws.Range("D17").Value = "0.000009139"
$ws.Range("E17").Value = "  +9.25%  "

# Row 18: Uniswap
$ws.Range("D18").Value = "5.980"
$ws.Range("E18").Value = "  -2.15%  "

# Row 19: WrappedBTC
$ws.Range("D19").Value = "29.040.29"
$ws.Range("E19").Value = "  -0.57%  "

# Row 20: WrappedliquidstakedEther2.0
$ws.Range("D20").Value = "2.075.76"
$ws.Range("E20").Value = "  -0.92%  "

# Row 21: BitcoinCash
$ws.Range("D21").Value = "224.69"
$ws.Range("E21").Value = "  -1.46%  "

# Row 22: Avalanche
$ws.Range("E22").Value = "  -1.18%  "

# Row 23: Dai
$ws.Range("E23").Value = "  -0.06%  "

# Row 24: Chainlink
$ws.Range("D24").Value = "7.174"
$ws.Range("E24").Value = "  +0.28%  "

# Row 25: BinanceUSD
$ws.Range("E25").Value = "  +0.04%  "

# Row 26: Monero
$ws.Range("D26").Value = "159.43"
$ws.Range("E26").Value = "  +0.07%  "

# Row 27: Cosmos
$ws.Range("D27").Value = "8.406"
$ws.Range("E27").Value = "  -2.48%  "

# Row 28: Stellar
$ws.Range("D28").Value = "0.1354"
$ws.Range("E28").Value = "  -4.12%  "

# Row 29: EthereumClassic
$ws.Range("D29").Value = "17.80"
$ws.Range("E29").Value = "  -1.03%  "

# Row 30: PancakeSwap
$ws.Range("E30").Value = "  -1.37%  "

# Row 31: InternetComputer(DFINITY)
$ws.Range("D31").Value = "4.028"
$ws.Range("E31").Value = "  -0.39%  "

# Row 32: Toncoin -> Filecoin
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "4.046"
$ws.Range("E32").Value = "  -1.60%  "

# Row 33: Filecoin -> Toncoin
$ws.Range("B33").Value = "Toncoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D33").Value = "1.205"
$ws.Range("E33").Value = "  +1.06%  "

# Row 34: Hedera
$ws.Range("D34").Value = "0.05214"
$ws.Range("E34").Value = "  -1.69%  "

# Row 35: LidoDAOToken
$ws.Range("D35").Value = "1.830"
$ws.Range("E35").Value = "  -1.87%  "

# Row 36: ARBITRUM
$ws.Range("D36").Value = "1.150"
$ws.Range("E36").Value = "  +0.98%  "

# Row 37: ImmutableX
$ws.Range("D37").Value = "0.7315"
$ws.Range("E37").Value = "  -1.97%  "

# Row 38: HuobiToken
$ws.Range("E38").Value = "  -0.20%  "

# Row 39: Maker
$ws.Range("D39").Value = "1.271.93"
$ws.Range("E39").Value = "  -1.12%  "

# Row 40: MXToken
$ws.Range("D40").Value = "2.747"
$ws.Range("E40").Value = "  +0.44%  "

# Row 41: VeChain
$ws.Range("D41").Value = "0.01780"
$ws.Range("E41").Value = "  -1.14%  "

# Row 42: FraxShare
$ws.Range("D42").Value = "6.303"
$ws.Range("E42").Value = "  +7.10%  "

# Row 43: TrustWalletToken
$ws.Range("D43").Value = "0.8933"
$ws.Range("E43").Value = "  -4.59%  "

# Row 44: PaxDollar
$ws.Range("E44").Value = "  +0.16%  "

# Row 45: Quant
$ws.Range("D45").Value = "101.86"
$ws.Range("E45").Value = "  -0.50%  "

# Row 46: RocketPoolETH
$ws.Range("D46").Value = "1.974.48"
$ws.Range("E46").Value = "  -0.80%  "

# Row 47: Mantle
$ws.Range("D47").Value = "0.5117"
$ws.Range("E47").Value = "  -0.57%  "

# Row 48: Aave
$ws.Range("D48").Value = "63.26"
$ws.Range("E48").Value = "  +0.33%  "

# Row 49: BabyDogeCoin
$ws.Range("E49").Value = "  -0.64%  "

# Row 50: TheSandbox
$ws.Range("D50").Value = "0.3958"
$ws.Range("E50").Value = "  -1.51%  "

# Row 51: RenderToken
$ws.Range("D51").Value = "1.671"
$ws.Range("E51").Value = "  -5.07%  "

# Restore default (unstyled) formatting on the touched price/volume range so
# the cells keep text storage without leaving a stray NumberFormat behind.
$ws.Range("D2:E51").Style = "Normal"
